$wb = $excel.ActiveWorkbook

# Sheet 1: 積み込み
$ws = $wb.Worksheets.Item(1)
$ws.Range("D1").Value = "←"
$ws.Range("F1").Value = "←"
$ws.Range("H1").Value = "←"
$ws.Range("J1").Value = "→"
$ws.Range("L1").Value = ""
$ws.Range("N1").Value = ""
$ws.Range("V1").Value = "→"
$ws.Range("Z1").Value = "←"
$ws.Range("AD1").Value = "→"
$ws.Range("C2").Value = "↓"
$ws.Range("I2").Value = "↑"
$ws.Range("K2").Value = ""
$ws.Range("M2").Value = "↑"
$ws.Range("O2").Value = ""
$ws.Range("W2").Value = ""
$ws.Range("Y2").Value = "↓"
$ws.Range("AC2").Value = "↓"
$ws.Range("AE2").Value = ""
$ws.Range("AH2").Value = "2 → 6"
$ws.Range("B3").Value = "←"
$ws.Range("F3").Value = ""
$ws.Range("H3").Value = "→"
$ws.Range("J3").Value = ""
$ws.Range("N3").Value = "→"
$ws.Range("R3").Value = "→"
$ws.Range("T3").Value = "→"
$ws.Range("V3").Value = "→"
$ws.Range("X3").Value = ""
$ws.Range("AD3").Value = "→"
$ws.Range("AH3").Value = "0 → 7"
$ws.Range("A4").Value = ""
$ws.Range("C4").Value = "↓"
$ws.Range("G4").Value = "↑"
$ws.Range("K4").Value = "↑"
$ws.Range("M4").Value = "↑"
$ws.Range("O4").Value = ""
$ws.Range("S4").Value = ""
$ws.Range("W4").Value = "↓"
$ws.Range("AC4").Value = ""
$ws.Range("AE4").Value = "↓"
$ws.Range("AH4").Value = "0 → 9"
$ws.Range("D5").Value = "→"
$ws.Range("F5").Value = ""
$ws.Range("P5").Value = "←"
$ws.Range("Z5").Value = "→"
$ws.Range("AB5").Value = ""
$ws.Range("AH5").Value = "0 → 7"
$ws.Range("C6").Value = ""
$ws.Range("I6").Value = ""
$ws.Range("O6").Value = ""
$ws.Range("U6").Value = "↓"
$ws.Range("W6").Value = ""
$ws.Range("Y6").Value = "↑"
$ws.Range("AA6").Value = "↓"
$ws.Range("AC6").Value = "↑"
$ws.Range("AE6").Value = ""
$ws.Range("AH6").Value = "2 → 9"
$ws.Range("B7").Value = "→"
$ws.Range("J7").Value = "←"
$ws.Range("L7").Value = "←"
$ws.Range("N7").Value = "←"
$ws.Range("P7").Value = "←"
$ws.Range("R7").Value = ""
$ws.Range("T7").Value = "←"
$ws.Range("Z7").Value = ""

# Sheet 2: 搬出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F1").Value = ""
$ws.Range("L1").Value = "→"
$ws.Range("Z1").Value = ""
$ws.Range("C2").Value = ""
$ws.Range("E2").Value = "↓"
$ws.Range("G2").Value = "↑"
$ws.Range("I2").Value = ""
$ws.Range("K2").Value = "↑"
$ws.Range("AA2").Value = "↓"
$ws.Range("AE2").Value = ""
$ws.Range("AH2").Value = "2 → 6"
$ws.Range("B3").Value = "←"
$ws.Range("F3").Value = ""
$ws.Range("J3").Value = "→"
$ws.Range("L3").Value = ""
$ws.Range("N3").Value = "→"
$ws.Range("R3").Value = ""
$ws.Range("V3").Value = "←"
$ws.Range("AD3").Value = "←"
$ws.Range("AH3").Value = "0 → 7"
$ws.Range("A4").Value = "↓"
$ws.Range("C4").Value = ""
$ws.Range("M4").Value = "↑"
$ws.Range("O4").Value = ""
$ws.Range("S4").Value = "↓"
$ws.Range("U4").Value = ""
$ws.Range("W4").Value = "↑"
$ws.Range("Y4").Value = "↑"
$ws.Range("AC4").Value = "↓"
$ws.Range("AE4").Value = ""
$ws.Range("AH4").Value = "0 → 9"
$ws.Range("B5").Value = ""
$ws.Range("D5").Value = ""
$ws.Range("F5").Value = ""
$ws.Range("J5").Value = "←"
$ws.Range("N5").Value = ""
$ws.Range("P5").Value = "→"
$ws.Range("R5").Value = "←"
$ws.Range("T5").Value = "←"
$ws.Range("V5").Value = ""
$ws.Range("AB5").Value = "←"
$ws.Range("AD5").Value = "←"
$ws.Range("AH5").Value = "0 → 7"
$ws.Range("A6").Value = "↓"
$ws.Range("C6").Value = "↓"
$ws.Range("E6").Value = "↓"
$ws.Range("G6").Value = "↓"
$ws.Range("I6").Value = "↓"
$ws.Range("K6").Value = ""
$ws.Range("M6").Value = ""
$ws.Range("O6").Value = "↑"
$ws.Range("U6").Value = ""
$ws.Range("W6").Value = "↑"
$ws.Range("Y6").Value = ""
$ws.Range("AA6").Value = "↑"
$ws.Range("AC6").Value = "↑"
$ws.Range("AH6").Value = "2 → 9"
$ws.Range("B7").Value = "→"
$ws.Range("D7").Value = "→"
$ws.Range("F7").Value = "→"
$ws.Range("H7").Value = "→"
$ws.Range("J7").Value = "→"
$ws.Range("L7").Value = "→"
$ws.Range("N7").Value = "→"
$ws.Range("T7").Value = "←"
$ws.Range("V7").Value = ""
$ws.Range("X7").Value = ""
$ws.Range("Z7").Value = "→"
$ws.Range("AB7").Value = ""

